$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G15").Value = "2025-05-18 23:35:01"
$ws.Range("G16").Value = "2025-05-18 23:35:01"
$ws.Range("G17").Value = "2025-05-18 23:35:01"
$ws.Range("G18").Value = "2025-05-18 23:35:02"
$ws.Range("G19").Value = "2025-05-18 23:35:02"
$ws.Range("G20").Value = "2025-05-18 23:35:02"
$ws.Range("G21").Value = "2025-05-18 23:35:02"
$ws.Range("G22").Value = "2025-05-18 23:35:03"
$ws.Range("G23").Value = "2025-05-18 23:35:03"
$ws.Range("G24").Value = "2025-05-18 23:35:03"
$ws.Range("G25").Value = "2025-05-18 23:35:03"
$ws.Range("G26").Value = "2025-05-18 23:35:03"
$ws.Range("G27").Value = "2025-05-18 23:35:04"
